$wb = $excel.ActiveWorkbook

# --- Settings sheet: add FilePath / SheetName / URL rows for the new
#     "InputData" config block (Initialization step update) ---
$settings = $wb.Worksheets.Item("Settings")

# Insert three new rows above the existing "OrchestratorQueueName" row,
# shifting all rows (and the trailing blank-formatted rows) down by 3.
$settings.Rows("2:4").Insert()

$settings.Range("A2").Value = "FilePath"
$settings.Range("B2").Value = "Data\Data.xlsx"

$settings.Range("A3").Value = "SheetName"
$settings.Range("B3").Value = "InputData"

$settings.Range("A4").Value = "URL"
$settings.Range("B4").Value = "https://rpachallenge.com/"

# Leave the selection on the newly added SheetName row, matching the
# author's saved cursor position.
$settings.Range("A3").Select()

# --- Constants sheet: "get transaction data" step now retries once
#     (MaxRetryNumber 0 -> 1) ---
$constants = $wb.Worksheets.Item("Constants")
$constants.Range("B2").Value = 1

# Constants becomes the active/visible tab with the cursor left on B2.
$constants.Activate()
$constants.Range("B2").Select()
